$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.088919401168823
$ws.Range("B1").Value = 2.936525106430054
$ws.Range("C1").Value = 2.2719886302948
$ws.Range("D1").Value = 2.133402109146118
$ws.Range("E1").Value = 2.130094289779663
